$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "oddballStimuli/trigger_Set2-Mult-Practice.wav"
$ws.Range("C2").Select() | Out-Null
